# AFDP-2327 - Add 'Fulfill' to the list of possible next queues for every
# rule whose "Default return queue" (column G) is "Fulfill" (i.e. rows 25-31
# of the "Next Possible Queues Rules" table on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = 25, 26, 27, 28, 29, 30, 31

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)   # column E = "List of possible next queues"
    $current = $cell.Value2
    if ($current -notlike "*Fulfill*") {
        $cell.Value2 = "$current,Fulfill"
    }
}
